$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.5306
$ws.Range("D3").Value = -7.013399999999993
$ws.Range("D14").Value = -7.887499999999994
$ws.Range("D16").Value = -8.763100000000003
$ws.Range("E18").Value = 18.67260000000001
$ws.Range("D21").Value = -8.429799999999995
$ws.Range("D23").Value = -7.052599999999995
$ws.Range("E24").Value = 16.4167
$ws.Range("D25").Value = -7.543599999999998
$ws.Range("E25").Value = 16.90630000000001
$ws.Range("D26").Value = -8.927399999999995
$ws.Range("E27").Value = 16.49099999999999
$ws.Range("D29").Value = -7.022199999999995
$ws.Range("E30").Value = 15.3037
$ws.Range("E31").Value = 15.8909
$ws.Range("E39").Value = 15.4372
$ws.Range("D40").Value = -8.806799999999988
$ws.Range("E42").Value = 16.6842
$ws.Range("E48").Value = 17.2992
$ws.Range("E51").Value = 17.1032
$ws.Range("E52").Value = 16.9858
$ws.Range("D53").Value = -6.346699999999995
$ws.Range("E55").Value = 16.592
$ws.Range("E56").Value = 15.3544
$ws.Range("D57").Value = -8.749999999999998
$ws.Range("E57").Value = 16.11620000000001
$ws.Range("D59").Value = -8.322200000000002
$ws.Range("E60").Value = 15.5552
$ws.Range("D65").Value = -7.906099999999999
$ws.Range("D69").Value = -7.077099999999997
$ws.Range("E73").Value = 17.88360000000001
$ws.Range("E74").Value = 16.6271
$ws.Range("D79").Value = -6.057999999999998
$ws.Range("D83").Value = -8.880599999999998
$ws.Range("E89").Value = 17.35640000000001
$ws.Range("E90").Value = 16.0063
$ws.Range("D91").Value = -6.060999999999996
$ws.Range("E92").Value = 18.65830000000002
$ws.Range("D93").Value = -6.510999999999995
$ws.Range("D100").Value = -8.310700000000004
